$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 617, shifting existing rows 617:658 down to 618:659
$ws.Rows.Item(617).Insert()

# Populate the newly inserted row 617 with the new record.
# Column A holds a date formatted as plain text (e.g. "2026/01/10"); setting
# NumberFormat to Text before assigning the value stops Excel from
# auto-converting the date-like string into a date serial number. The
# format is then reset back to Normal so no stray style is left behind.
$ws.Cells.Item(617, 1).NumberFormat = "@"
$ws.Cells.Item(617, 1).Value = "2026/01/10"
$ws.Cells.Item(617, 1).Style = "Normal"

$ws.Cells.Item(617, 2).Value = "土"
$ws.Cells.Item(617, 3).Value = 6
$ws.Cells.Item(617, 4).Value = 201
